# PfsIcsMhsHardwareConfig.pptx
# Commit: "removed xxx machine type from host#2 description"
#
# The "Host #2" box on slide 1 described the machine as
# "(XXX R710 running Debian 9.x)". The placeholder machine-type prefix
# ("XXX R710 ") is removed, leaving "(Running Debian 9.x)".
#
# We locate the run by its old text and rewrite only the exact
# substring via TextRange.Characters(start,length) so the surrounding
# run/paragraph formatting (size, color, etc.) is left completely
# untouched - matching the original run's <a:rPr> in the OOXML.

$p = $ppt.ActivePresentation

$oldFragment = "(XXX R710 running Debian 9.x)"
$newFragment = "(Running Debian 9.x)"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            $fullText = $tr.Text

            $idx = $fullText.IndexOf($oldFragment)
            if ($idx -ge 0) {
                # TextRange.Characters is 1-based.
                $start = $idx + 1
                $length = $oldFragment.Length

                $sub = $tr.Characters($start, $length)
                $sub.Text = $newFragment

                Write-Host "Updated slide $si shape $shi ($($shape.Name)): '$oldFragment' -> '$newFragment'"
            }
        }
    }
}
